{"js": "// Update the \"vX.Y (YYYY-MM-DD)\" version/date line on the title page.\n// The paragraph (style \"Date\") holds three runs: version, a space, and the\n// parenthesized date. We replace the version text and the date text,\n// leaving the middle space run untouched.\n\nconst body = context.document.body;\n\nconst versionResults = body.search(\"v1.2\", { matchCase: true, matchWholeWord: false });\nversionResults.load(\"items,text\");\nawait context.sync();\n\nif (versionResults.items.length > 0) {\n  versionResults.items[0].insertText(\"v1.3\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst dateResults = body.search(\"(2014-05-28)\", { matchCase: true, matchWholeWord: false });\ndateResults.load(\"items,text\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"(2014-06-03)\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the \"vX.Y (YYYY-MM-DD)\" version/date line on the title page.\n# The paragraph (style \"Date\") holds three runs: version, a space, and the\n# parenthesized date. We replace the version text and the date text via\n# Find & Replace, leaving the untouched space run alone.\n\n$d = $word.ActiveDocument\n\n$wdReplaceOne   = 1\n$wdFindContinue = 1\n\n# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n#              MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n\n$range1 = $d.Content\n$range1.Find.ClearFormatting()\n$range1.Find.Execute(\"v1.2\", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"v1.3\", $wdReplaceOne)\n\n$range2 = $d.Content\n$range2.Find.ClearFormatting()\n$range2.Find.Execute(\"(2014-05-28)\", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"(2014-06-03)\", $wdReplaceOne)\n"}
